$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update memory-layout formulas
$ws.Range("G2").Formula = '=$B$2+256+256'
$ws.Range("H2").Formula = '=$B$2+256'
$ws.Range("C3").Formula = '=26*1024'

# New row 5 labels for the task/queue columns
$ws.Range("G5").Value = "irrigation"
$ws.Range("H5").Value = "wireless"
$ws.Range("I5").Value = "adc"
$ws.Range("J5").Value = "status"
$ws.Range("K5").Value = "user button"

# Update the active selection to match the authored state
$ws.Range("H4").Select()
